$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("labels")
$ws.Activate()

$ws.Range("A33").Value = "measuresOfCenter"
$ws.Range("B33").Value = "subsection"
$ws.Range("C33").Value = "Note change from OI, careful about future refs"

$ws.Range("A34").Value = "measuresOfSpread"
$ws.Range("B34").Value = "subsection"
$ws.Range("C34").Value = "Note change from OI"

$ws.Range("A32:C32").Copy()
$ws.Range("A33:C34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C35").Select()
